$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 gets the first new string value
$ws.Range("F2").Value = "lksjdjflkj"

# G6 gets the second new string value, and ends up as the active selection
$ws.Range("G6").Value = "kjsldkjkf"
$ws.Range("G6").Select()
